# Remove the "[Work]" / "(803) 799-0425" row from the Assets table.
#
# The document has two Field/Value tables. The second one (the attorney
# info table) has a row for [Work] with the phone number (803) 799-0425.
# That whole table row is being removed.

$d = $word.ActiveDocument

foreach ($table in $d.Tables) {
    for ($i = $table.Rows.Count; $i -ge 1; $i--) {
        $row = $table.Rows.Item($i)
        $text = $row.Range.Text
        if ($text.Contains("[Work]") -and $text.Contains("(803) 799-0425")) {
            $row.Delete()
        }
    }
}
